# Add Test Data for Russia / Finland / Hungary markets.
# Each new sheet is a copy of the last existing sheet ("Denmark"), renamed,
# with the "Description" (B2) and "User Story" (B4) cells updated.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$template = $wb.Sheets.Item("Denmark")

# Burn one internal sheetId (19) on a throwaway sheet so the three new
# sheets line up on sheetId 20/21/22 (matching upstream history), then
# discard it without leaving a trace in the tab order. It must stay alive
# until after the real copies are made, otherwise the id gets reused.
$placeholder = $wb.Sheets.Add([Type]::Missing, $template)

# --- Russia -----------------------------------------------------------
$template.Copy([Type]::Missing, $placeholder)
$russia = $wb.Sheets.Item($placeholder.Index + 1)
$russia.Name = "Russia"
$russia.Range("B2").Value = "Russia Market"
$russia.Range("B4").Value = "NGC-2929/T2897"
$russia.Activate()
$russia.Range("A1:D17").Select() | Out-Null

# --- Finland ------------------------------------------------------------
$russia.Copy([Type]::Missing, $russia)
$finland = $wb.Sheets.Item($russia.Index + 1)
$finland.Name = "Finland"
$finland.Range("B2").Value = "Finland Market"
$finland.Range("B4").Value = "NGC-3130/T2940"
$finland.Activate()
$finland.Range("A1:D17").Select() | Out-Null

# --- Hungary --------------------------------------------------------------
$finland.Copy([Type]::Missing, $finland)
$hungary = $wb.Sheets.Item($finland.Index + 1)
$hungary.Name = "Hungary"
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Range("B4").Value = "NGC-3104/T2989"
$hungary.Activate()
$hungary.Range("A1:D17").Select() | Out-Null

$placeholder.Delete() | Out-Null

# Re-fetch by name: deleting a sheet earlier in the tab order can leave a
# previously-captured reference pointing at a stale Index.
$hungary = $wb.Sheets.Item("Hungary")
$hungary.Activate()
